$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.346238017082214
$ws.Range("B1").Value = 1.447638392448425
$ws.Range("C1").Value = 1.547587037086487
$ws.Range("D1").Value = 2.177873134613037
$ws.Range("E1").Value = 3.764452695846558
